# This fix addresses several issues with updating CSV files when a form has a repository.
# Concretely: insert a new "Project" question (select_one_from_file generated_case.csv)
# into the "survey" sheet of this XLSForm workbook, right before the "coll_date" (date) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new row above row 5 ("date" / coll_date), pushing the rest of the
# questions down by one. Excel copies the formatting of the row above (row 4)
# into the freshly inserted row.
$ws.Rows(5).Insert()

# Populate the new row with the "Project" question.
$ws.Cells.Item(5, 1).Value = "select_one_from_file generated_case.csv"
$ws.Cells.Item(5, 2).Value = "me_project"
$ws.Cells.Item(5, 3).Value = "Project"

# Leave the new row's "required"/"choice_filter" columns blank (already
# empty after the insert/copy-down).

# Reflect the user's last selection in the sheet (cell B5, the new "name" cell).
$ws.Range("B5").Select() | Out-Null
